# Insert 4 new rows at 902-905 (shifts existing rows 902-974 down to 906-978)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("902:905").Insert()

# Common (constant) column values for this block of rows
$marketId = 7
$market   = "Terminal Hortofrutícola Agro Chillán"
$region   = "Ñuble"
$codreg   = 16
$tipo     = "Fruta"
$prodId   = 100104
$prod     = "Frutos de pepita"
$catId    = 100104002
$cat      = "Manzana"
$unidad   = "$/caja 16 kilos empedrada"
$origen   = "Provincia de Curicó"
$kgUnit   = 16

$fecha = 44826

# Row 902: Fuji royal, Especial
$r = 902
$ws.Cells.Item($r, 1).Value  = $marketId
$ws.Cells.Item($r, 2).Value  = $market
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $prodId
$ws.Cells.Item($r, 8).Value  = $prod
$ws.Cells.Item($r, 9).Value  = $catId
$ws.Cells.Item($r, 10).Value = $cat
$ws.Cells.Item($r, 11).Value = "Fuji royal"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 60
$ws.Cells.Item($r, 14).Value = 11000
$ws.Cells.Item($r, 15).Value = 11000
$ws.Cells.Item($r, 16).Value = 11000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 688
$ws.Cells.Item($r, 20).Value = $kgUnit

# Row 903: Fuji royal, Primera
$r = 903
$ws.Cells.Item($r, 1).Value  = $marketId
$ws.Cells.Item($r, 2).Value  = $market
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $prodId
$ws.Cells.Item($r, 8).Value  = $prod
$ws.Cells.Item($r, 9).Value  = $catId
$ws.Cells.Item($r, 10).Value = $cat
$ws.Cells.Item($r, 11).Value = "Fuji royal"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 120
$ws.Cells.Item($r, 14).Value = 9500
$ws.Cells.Item($r, 15).Value = 10000
$ws.Cells.Item($r, 16).Value = 9750
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 609
$ws.Cells.Item($r, 20).Value = $kgUnit

# Row 904: Granny Smith, Especial
$r = 904
$ws.Cells.Item($r, 1).Value  = $marketId
$ws.Cells.Item($r, 2).Value  = $market
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $prodId
$ws.Cells.Item($r, 8).Value  = $prod
$ws.Cells.Item($r, 9).Value  = $catId
$ws.Cells.Item($r, 10).Value = $cat
$ws.Cells.Item($r, 11).Value = "Granny Smith"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 60
$ws.Cells.Item($r, 14).Value = 12000
$ws.Cells.Item($r, 15).Value = 12000
$ws.Cells.Item($r, 16).Value = 12000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 750
$ws.Cells.Item($r, 20).Value = $kgUnit

# Row 905: Granny Smith, Primera
$r = 905
$ws.Cells.Item($r, 1).Value  = $marketId
$ws.Cells.Item($r, 2).Value  = $market
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $prodId
$ws.Cells.Item($r, 8).Value  = $prod
$ws.Cells.Item($r, 9).Value  = $catId
$ws.Cells.Item($r, 10).Value = $cat
$ws.Cells.Item($r, 11).Value = "Granny Smith"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 120
$ws.Cells.Item($r, 14).Value = 10000
$ws.Cells.Item($r, 15).Value = 11000
$ws.Cells.Item($r, 16).Value = 10500
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 656
$ws.Cells.Item($r, 20).Value = $kgUnit
